$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = [double]"0"
$ws.Range("E2").Value = [double]"23.8900000000003"
$ws.Range("G2").Value = [double]"0.0002881931882366962"
$ws.Range("H2").Value = [double]"0.003682274204135995"
$ws.Range("K2").Value = [double]"4.619230259713326"
$ws.Range("L2").Value = "[1.8369814559780275, 7.401479063448625]"
$ws.Range("M2").Value = [double]"0.001196164188009829"
$ws.Range("N2").Value = [double]"0.002392328376019659"
$ws.Range("O2").Value = [double]"0.2830263651882685"
$ws.Range("P2").Value = "[-0.42139481039142623, 0.9874475407679633]"
$ws.Range("Q2").Value = [double]"0.4300220547650238"
$ws.Range("R2").Value = [double]"0.4300220547650238"
$ws.Range("S2").Value = [double]"13.75867878657695"
$ws.Range("T2").Value = "[12.183508275256614, 15.333849297897281]"
$ws.Range("W2").Value = [double]"22.81387387387416"
$ws.Range("X2").Value = [double]"20.13551551551576"
$ws.Range("Y2").Value = [double]"25.49223223223256"

# Row 3
$ws.Range("E3").Value = [double]"24.5400000000004"
$ws.Range("G3").Value = [double]"0.006313871695905693"
$ws.Range("H3").Value = [double]"0.01622267331555555"
$ws.Range("K3").Value = [double]"4.459126342152437"
$ws.Range("L3").Value = "[0.8933080048093451, 8.024944679495528]"
$ws.Range("M3").Value = [double]"0.01443107058536386"
$ws.Range("N3").Value = [double]"0.01443107058536386"
$ws.Range("O3").Value = [double]"2.496921488438735"
$ws.Range("P3").Value = "[1.7421845146033483, 3.251658462274121]"
$ws.Range("Q3").Value = [double]"3.343352261708787e-10"
$ws.Range("R3").Value = [double]"6.686704523417575e-10"
$ws.Range("S3").Value = [double]"14.25723220270636"
$ws.Range("T3").Value = "[12.319546565352947, 16.194917840059777]"
$ws.Range("W3").Value = [double]"14.78786786786811"
$ws.Range("X3").Value = [double]"11.84012012012031"
$ws.Range("Y3").Value = [double]"17.73561561561591"

$wb.Save()
